$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 87; this shifts old rows 87:89 down to 89:91
$ws.Rows(87).Insert()
$ws.Rows(88).Insert()

# New row 87: Primera, date 2021-09-09 (44448)
$ws.Range("A87").Value = 11
$ws.Range("B87").Value = "Vega Monumental Concepción"
$ws.Range("C87").Value = "Bíobío"
$ws.Range("D87").Value = 44448
$ws.Range("E87").Value = 8
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100101
$ws.Range("H87").Value = "Berries"
$ws.Range("I87").Value = 100101007
$ws.Range("J87").Value = "Kiwi"
$ws.Range("K87").Value = "Hayward"
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 200
$ws.Range("N87").Value = 13000
$ws.Range("O87").Value = 14000
$ws.Range("P87").Value = 13500
$ws.Range("Q87").Value = "$/bandeja 18 kilos"
$ws.Range("R87").Value = "Región de O'Higgins"
$ws.Range("S87").Value = 750
$ws.Range("T87").Value = 18

# New row 88: Segunda, date 2021-09-09 (44448)
$ws.Range("A88").Value = 11
$ws.Range("B88").Value = "Vega Monumental Concepción"
$ws.Range("C88").Value = "Bíobío"
$ws.Range("D88").Value = 44448
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = "Fruta"
$ws.Range("G88").Value = 100101
$ws.Range("H88").Value = "Berries"
$ws.Range("I88").Value = 100101007
$ws.Range("J88").Value = "Kiwi"
$ws.Range("K88").Value = "Hayward"
$ws.Range("L88").Value = "Segunda"
$ws.Range("M88").Value = 100
$ws.Range("N88").Value = 12000
$ws.Range("O88").Value = 12000
$ws.Range("P88").Value = 12000
$ws.Range("Q88").Value = "$/bandeja 18 kilos"
$ws.Range("R88").Value = "Región de O'Higgins"
$ws.Range("S88").Value = 667
$ws.Range("T88").Value = 18

# Apply the date number format (used by existing D column cells) to the new D cells
$ws.Range("D87").NumberFormat = $ws.Range("D86").NumberFormat
$ws.Range("D88").NumberFormat = $ws.Range("D86").NumberFormat

$wb.Save()
